$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 426; existing rows 426-486 shift down to 427-487.
$ws.Rows("426:426").Insert()

# Populate the newly inserted row 426 with the new data record.
$ws.Range("A426").Value = 1
$ws.Range("B426").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C426").Value = "Arica y Parinacota"
$ws.Range("D426").Value = 45154
$ws.Range("E426").Value = 15
$ws.Range("F426").Value = "Fruta"
$ws.Range("G426").Value = 100102
$ws.Range("H426").Value = "Cítricos"
$ws.Range("I426").Value = 100102003
$ws.Range("J426").Value = "Limón"
$ws.Range("K426").Value = "Tahití"
$ws.Range("L426").Value = "Primera"
$ws.Range("M426").Value = 400
$ws.Range("N426").Value = 43000
$ws.Range("O426").Value = 45000
$ws.Range("P426").Value = 44250
$ws.Range("Q426").Value = "`$/caja 24 kilos"
$ws.Range("R426").Value = "Perú"
$ws.Range("S426").Value = 1844
$ws.Range("T426").Value = 24
